$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("F 0.3")
$dst = $wb.Worksheets.Item("F 0.5")

# Select and copy the source range on "F 0.3" (A2:G8)
$src.Activate()
$src.Range("A2:G8").Select()
$src.Range("A2:G8").Copy()

# Paste into "F 0.5" starting at A13
$dst.Activate()
$dst.Range("A13").Select()
$dst.Paste()

# Restore the number formats that the source cells carried
# (D17:D19 show as dates/2-decimal, B18 as an integer in the source range)
$dst.Range("D17:D19").NumberFormat = "0.00"
$dst.Range("B18").NumberFormat = "0"

# Set final selection on "F 0.5"
$dst.Range("F8").Select()
